# Created bat file to fetch missing info
# -> appends the newly-fetched rows to the "repair_notes" and "new_inventory"
#    worksheets.

$wb = $excel.ActiveWorkbook

# Helper: write a value as literal text so number-/date-looking strings
# (e.g. "10-15-2020", "5555", "1") are not auto-coerced into numbers/dates.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- repair_notes: new log entry (row 11) -----------------------------
$ws1 = $wb.Worksheets.Item("repair_notes")

Set-TextValue $ws1.Range("A11") "8264CLM"
Set-TextValue $ws1.Range("B11") "Windows-10-10.0.19041-SP0"
Set-TextValue $ws1.Range("C11") "Intel(R) Core(TM) i5-7300U CPU @ 2.60GHz"
Set-TextValue $ws1.Range("D11") "8 GB"
Set-TextValue $ws1.Range("E11") "10.110.200.112"
Set-TextValue $ws1.Range("F11") "10-15-2020"
Set-TextValue $ws1.Range("G11") "09:37"
Set-TextValue $ws1.Range("H11") "Home"
Set-TextValue $ws1.Range("I11") "My PC"
Set-TextValue $ws1.Range("J11") "Dan"

# --- new_inventory: new inventory row (row 10) -------------------------
$ws2 = $wb.Worksheets.Item("new_inventory")

Set-TextValue $ws2.Range("A10") "8264CLM"
Set-TextValue $ws2.Range("B10") "Y"
Set-TextValue $ws2.Range("C10") "Dell"
Set-TextValue $ws2.Range("D10") "XPS 15"
Set-TextValue $ws2.Range("E10") "test"
Set-TextValue $ws2.Range("F10") "5555"
Set-TextValue $ws2.Range("H10") "Dan"
Set-TextValue $ws2.Range("I10") "Home"
Set-TextValue $ws2.Range("Q10") "1"
Set-TextValue $ws2.Range("U10") "1"
Set-TextValue $ws2.Range("X10") "1"
Set-TextValue $ws2.Range("Y10") "1"

Write-Host "Appended new rows to repair_notes and new_inventory"
